# Auto-generated Excel COM-interop script
# Updates the cryptos list: refreshed prices and 1h volume percentages,
# and two coin rows changed (ImmutableX <-> InternetComputer(DFINITY) swap
# order/values, dogwifhat replaced by Monero), matching the commit
# 'Updated cryptos list ... with GitHub Actions'.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.204.13'
$ws.Range('E2').Value = '  -6.24%  '

$ws.Range('D3').Value = '3.266.85'
$ws.Range('E3').Value = '  -7.85%  '

$ws.Range('E4').Value = '  +0.25%  '

$c = $ws.Range('D5')
$c.Value = '''179.59'
$c.Style = "Normal"
$ws.Range('E5').Value = '  -11.72%  '

$c = $ws.Range('D6')
$c.Value = '''515.79'
$c.Style = "Normal"
$ws.Range('E6').Value = '  -6.95%  '

$c = $ws.Range('D7')
$c.Value = '''0.594'
$c.Style = "Normal"
$ws.Range('E7').Value = '  -1.20%  '

$ws.Range('D8').Value = '3.258.76'
$ws.Range('E8').Value = '  -7.58%  '

$ws.Range('E9').Value = '  +0.08%  '

$c = $ws.Range('D10')
$c.Value = '''0.613'
$c.Style = "Normal"
$ws.Range('E10').Value = '  -7.15%  '

$c = $ws.Range('D11')
$c.Value = '''56.94'
$c.Style = "Normal"
$ws.Range('E11').Value = '  -7.04%  '

$ws.Range('E12').Value = '  -9.83%  '

$c = $ws.Range('D13')
$c.Value = '''0.0000253'
$c.Style = "Normal"
$ws.Range('E13').Value = '  -7.32%  '

$c = $ws.Range('D14')
$c.Value = '''9.03'
$c.Style = "Normal"
$ws.Range('E14').Value = '  -9.09%  '

$ws.Range('D15').Value = '3.799.90'
$ws.Range('E15').Value = '  -7.56%  '

$ws.Range('E16').Value = '  -6.19%  '

$ws.Range('D17').Value = '3.278.04'
$ws.Range('E17').Value = '  -7.51%  '

$c = $ws.Range('D18')
$c.Value = '''17.58'
$c.Style = "Normal"
$ws.Range('E18').Value = '  -6.20%  '

$ws.Range('D19').Value = '63.260.27'
$ws.Range('E19').Value = '  -5.94%  '

$c = $ws.Range('D20')
$c.Value = '''10.85'
$c.Style = "Normal"
$ws.Range('E20').Value = '  -9.31%  '

$c = $ws.Range('D21')
$c.Value = '''0.942'
$c.Style = "Normal"
$ws.Range('E21').Value = '  -9.76%  '

$c = $ws.Range('D22')
$c.Value = '''368.46'
$c.Style = "Normal"
$ws.Range('E22').Value = '  -5.96%  '

$c = $ws.Range('D23')
$c.Value = '''11.17'
$c.Style = "Normal"
$ws.Range('E23').Value = '  -8.43%  '

$c = $ws.Range('D24')
$c.Value = '''3.67'
$c.Style = "Normal"
$ws.Range('E24').Value = '  -9.66%  '

$c = $ws.Range('D25')
$c.Value = '''79.58'
$c.Style = "Normal"
$ws.Range('E25').Value = '  -4.19%  '

$c = $ws.Range('D26')
$c.Value = '''3.80'
$c.Style = "Normal"
$ws.Range('E26').Value = '  -0.07%  '

$ws.Range('E27').Value = '  -1.95%  '

$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range('D28')
$c.Value = '''11.30'
$c.Style = "Normal"
$ws.Range('E28').Value = '  -6.85%  '

$ws.Range('B29').Value = 'ImmutableX'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range('D29')
$c.Value = '''2.60'
$c.Style = "Normal"
$ws.Range('E29').Value = '  -7.80%  '

$c = $ws.Range('D30')
$c.Value = '''8.27'
$c.Style = "Normal"
$ws.Range('E30').Value = '  -7.54%  '

$c = $ws.Range('D31')
$c.Value = '''28.28'
$c.Style = "Normal"
$ws.Range('E31').Value = '  -8.57%  '

$c = $ws.Range('D32')
$c.Value = '''634.83'
$c.Style = "Normal"
$ws.Range('E32').Value = '  -10.77%  '

$c = $ws.Range('D33')
$c.Value = '''6.63'
$c.Style = "Normal"
$ws.Range('E33').Value = '  -9.78%  '

$c = $ws.Range('D34')
$c.Value = '''11.10'
$c.Style = "Normal"
$ws.Range('E34').Value = '  -6.24%  '

$c = $ws.Range('D35')
$c.Value = '''0.105'
$c.Style = "Normal"
$ws.Range('E35').Value = '  -5.52%  '

$c = $ws.Range('D36')
$c.Value = '''58.35'
$c.Style = "Normal"
$ws.Range('E36').Value = '  -7.87%  '

$ws.Range('E37').Value = '  -0.04%  '

$c = $ws.Range('D38')
$c.Value = '''0.389'
$c.Style = "Normal"
$ws.Range('E38').Value = '  -5.41%  '

$c = $ws.Range('D39')
$c.Value = '''35.88'
$c.Style = "Normal"
$ws.Range('E39').Value = '  -11.54%  '

$c = $ws.Range('D40')
$c.Value = '''1.00'
$c.Style = "Normal"
$ws.Range('E40').Value = '  +0.19%  '

$ws.Range('D41').Value = '2.957.08'
$ws.Range('E41').Value = '  -5.40%  '

$c = $ws.Range('D42')
$c.Value = '''0.123'
$c.Style = "Normal"
$ws.Range('E42').Value = '  -4.81%  '

$ws.Range('D43').Value = '0.0₃0646'
$ws.Range('E43').Value = '  -9.52%  '

$c = $ws.Range('D44')
$c.Value = '''2.44'
$c.Style = "Normal"
$ws.Range('E44').Value = '  -4.19%  '

$c = $ws.Range('D45')
$c.Value = '''2.65'
$c.Style = "Normal"
$ws.Range('E45').Value = '  -14.82%  '

$c = $ws.Range('D46')
$c.Value = '''2.60'
$c.Style = "Normal"
$ws.Range('E46').Value = '  -4.58%  '

$ws.Range('E47').Value = '  -4.04%  '

$c = $ws.Range('D48')
$c.Value = '''2.78'
$c.Style = "Normal"

$c = $ws.Range('D49')
$c.Value = '''0.124'
$c.Style = "Normal"
$ws.Range('E49').Value = '  -2.91%  '

$c = $ws.Range('D50')
$c.Value = '''2.95'
$c.Style = "Normal"
$ws.Range('E50').Value = '  -3.02%  '

$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range('D51')
$c.Value = '''130.15'
$c.Style = "Normal"
$ws.Range('E51').Value = '  -5.14%  '

